$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "28.011.16"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").Value = "1.868.68"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.05"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5151"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3835"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08279"
$ws.Range("E9").Value = "  -3.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.110"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.54"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.198"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").Value = "1.871.41"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.308"
$ws.Range("E15").Value = "  +1.41%  "

$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.77"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06642"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.032"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").Value = "28.050.88"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("E24").Value = "  -2.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.253"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("D26").Value = "2.080.19"
$ws.Range("E26").Value = "  -0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.505"
$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.51"
$ws.Range("E28").Value = "  +0.37%  "

$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.030"
$ws.Range("E32").Value = "  -2.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.787"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.593"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.473"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06514"
$ws.Range("E37").Value = "  -1.14%  "

$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6598"
$ws.Range("E39").Value = "  +3.36%  "

$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.007"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.207"
$ws.Range("E42").Value = "  -3.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.19"
$ws.Range("E43").Value = "  -2.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6137"
$ws.Range("E44").Value = "  +2.21%  "

$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.672"
$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("E48").Value = "  +1.72%  "

$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.90"
$ws.Range("E50").Value = "  -0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.40"
$ws.Range("E51").Value = "  -2.81%  "
